# The deck's slide master and every one of its custom (slide) layouts carry
# a "Date Placeholder" (ppPlaceholderDate = 16) whose text is the cached
# value of an auto-updating datetimeFigureOut field. The commit that
# regenerated the golden file re-saved the reference deck on a later day,
# which re-stamped that cached value from 2/17/2018 to 4/5/2019 everywhere
# it appears (the slide master plus all 11 layouts). Reproduce that here.

$newDate = "4/5/2019"

$p = $ppt.ActivePresentation
$m = $p.SlideMaster

# Slide master's own Date Placeholder.
for ($j = 1; $j -le $m.Shapes.Count; $j++) {
    $sh = $m.Shapes.Item($j)
    if ($sh.PlaceholderFormat.Type -eq 16) {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

# Every custom layout hanging off the master.
$cl = $m.CustomLayouts
for ($i = 1; $i -le $cl.Count; $i++) {
    $lay = $cl.Item($i)
    for ($j = 1; $j -le $lay.Shapes.Count; $j++) {
        $sh = $lay.Shapes.Item($j)
        if ($sh.PlaceholderFormat.Type -eq 16) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}
